$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.678.88'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '1.725.82'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4929'
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06210'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").Value = '1.728.99'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.491'
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.15'
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9982'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").Value = '26.495.87'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9978'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007221'
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("D21").Value = '1.952.41'
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.461'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.531'
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '137.64'
$ws.Range("E28").Value = '  -1.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.07'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.914'
$ws.Range("E30").Value = '  -1.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.678'
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9360'
$ws.Range("E37").Value = '  +2.97%  '
$ws.Range("E38").Value = '  +3.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.411'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.602'
$ws.Range("E42").Value = '  +3.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.38'
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.868'
$ws.Range("E45").Value = '  +2.80%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.831'
$ws.Range("E48").Value = '  +1.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.56'
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("E51").Value = '  -1.93%  '
